# New crime data collected — weekly CompStat update for the 116th Precinct.
# Bumps the report "Volume ... Number" and the "Week Covering" date range by
# one week, and refreshes the Crime Complaints table (rows 15-31) with the
# newly collected weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text: volume number + reporting week dates -------------------
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# ---- Helper: write a numeric value while preserving/forcing a specific ----
# ---- numeric style, reusing an existing donor cell's format so no new  ----
# ---- cell style gets minted in styles.xml.                             ----
function Set-NumValue {
    param($range, $value, $donor)
    $range.Value = $value
    if ($donor) {
        $donor.Copy()
        $range.PasteSpecial(-4122)
    }
}

function Set-TextValue {
    param($range, $text, $donor)
    $range.Value = "'" + $text
    if ($donor) {
        $donor.Copy()
        $range.PasteSpecial(-4122)
    }
}

# Donor cells carrying the two numeric styles (s=14 plain count, s=15 pct)
# and the text style (s=13) used throughout the crime-complaints table.
$donorNum = $ws.Range("F15")     # style s=14
$donorPct = $ws.Range("H15")     # style s=15 (before this script edits it)
$donorTxt = $ws.Range("D15")     # style s=13 (text)

# Row 15 - Rape
Set-NumValue $ws.Range("C15") 2 $donorNum
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 16
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 220

# Row 16 - Robbery
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 24
$ws.Range("J16").Value = 32
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = -22.580645161290

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 93
$ws.Range("K17").Value = 7.526881720430
$ws.Range("L17").Value = 5.263157894736

# Row 18 - Burglary
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -50
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -37.931034482758
$ws.Range("L18").Value = -37.931034482758

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -56.521739130434
$ws.Range("I19").Value = 89
$ws.Range("J19").Value = 113
$ws.Range("K19").Value = -21.238938053097
$ws.Range("L19").Value = -19.090909090909

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 78
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = 25.806451612903
$ws.Range("L20").Value = 4

# Row 21 - TOTAL (bold styles s=17/18/19, already correct - values only)
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 21.428571428571
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -16.216216216216
$ws.Range("I21").Value = 325
$ws.Range("J21").Value = 335
$ws.Range("K21").Value = -2.985074626865
$ws.Range("L21").Value = -5.797101449275

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -55.555555555555
$ws.Range("F24").Value = 29
$ws.Range("H24").Value = -36.956521739130
$ws.Range("I24").Value = 191
$ws.Range("J24").Value = 242
$ws.Range("K24").Value = -21.074380165289
$ws.Range("L24").Value = -4.975124378109

# Row 25 - Retail Theft
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = -41.666666666666
$ws.Range("I25").Value = 44
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = -12
$ws.Range("L25").Value = 22.222222222222

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 50
$ws.Range("H26").Value = 47.058823529411
$ws.Range("I26").Value = 165
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 3.125
$ws.Range("L26").Value = 35.245901639344

# Row 27 - UCR Rape* (C27 text "0" -> numeric 2)
Set-NumValue $ws.Range("C27") 2 $donorNum
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = 171.428571428571
$ws.Range("L27").Value = 72.727272727272

# Row 28 - Other Sex Crimes (C28 numeric -> text "0"; D28 text "0" -> numeric 1;
# E28 text "***.*" -> numeric -100)
Set-TextValue $ws.Range("C28") "0" $donorTxt
Set-NumValue $ws.Range("D28") 1 $donorNum
Set-NumValue $ws.Range("E28") -100 $donorPct
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 200
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = -28.571428571428
$ws.Range("L28").Value = -64.285714285714

# Row 31 - Hate Crimes (D31 numeric 1 -> text "0"; E31 numeric -100 -> text "***.*")
Set-TextValue $ws.Range("D31") "0" $donorTxt
Set-TextValue $ws.Range("E31") "***.*" $donorTxt

Write-Output "CompStat weekly figures updated"
